# Updated cryptos list (Price / Volume(1h) columns) per the scraped diff.
# D-column values that look like a plain decimal number (e.g. "33.81") are
# written with a leading apostrophe so Excel stores them as text (quote-
# prefixed), matching the original inline-string cells instead of letting
# COM auto-convert them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.801.81'
$ws.Range('E2').Value = '  -1.15%  '
$ws.Range('D3').Value = '2.788.52'
$ws.Range('E3').Value = '  -2.47%  '
$ws.Range('D5').Value = '''358.56'
$ws.Range('E5').Value = '  -1.76%  '
$ws.Range('D6').Value = '''109.01'
$ws.Range('E6').Value = '  -6.32%  '
$ws.Range('D7').Value = '''0.557'
$ws.Range('E7').Value = '  +0.68%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '''0.588'
$ws.Range('E9').Value = '  -3.38%  '
$ws.Range('E10').Value = '  -6.61%  '
$ws.Range('D11').Value = '''0.0846'
$ws.Range('E11').Value = '  -2.47%  '
$ws.Range('E12').Value = '  +0.21%  '
$ws.Range('D13').Value = '''19.41'
$ws.Range('E13').Value = '  -3.73%  '
$ws.Range('D14').Value = '''7.62'
$ws.Range('E14').Value = '  -3.56%  '
$ws.Range('D15').Value = '3.232.21'
$ws.Range('E15').Value = '  -2.33%  '
$ws.Range('D16').Value = '2.797.87'
$ws.Range('E16').Value = '  -2.31%  '
$ws.Range('D17').Value = '''0.910'
$ws.Range('E17').Value = '  +0.18%  '
$ws.Range('D18').Value = '51.656.02'
$ws.Range('E18').Value = '  -1.62%  '
$ws.Range('D19').Value = '''7.37'
$ws.Range('E19').Value = '  +0.38%  '
$ws.Range('E20').Value = '  -2.94%  '
$ws.Range('D21').Value = '''13.05'
$ws.Range('E21').Value = '  -6.07%  '
$ws.Range('D22').Value = '0.0₃0975'
$ws.Range('E22').Value = '  -1.61%  '
$ws.Range('D23').Value = '''272.36'
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').Value = '''69.57'
$ws.Range('E24').Value = '  -1.58%  '
$ws.Range('E25').Value = '  -3.42%  '
$ws.Range('D26').Value = '''26.42'
$ws.Range('E26').Value = '  -2.93%  '
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('D28').Value = '''10.09'
$ws.Range('E28').Value = '  -2.50%  '
$ws.Range('E29').Value = '  -1.41%  '
$ws.Range('E30').Value = '  -1.00%  '
$ws.Range('D31').Value = '''0.0464'
$ws.Range('E31').Value = '  +2.63%  '
$ws.Range('D32').Value = '''51.36'
$ws.Range('E32').Value = '  +0.17%  '
$ws.Range('D33').Value = '''33.81'
$ws.Range('E33').Value = '  -2.38%  '
$ws.Range('D34').Value = '''5.72'
$ws.Range('E34').Value = '  -2.53%  '
$ws.Range('D35').Value = '''5.37'
$ws.Range('E35').Value = '  +7.07%  '
$ws.Range('D36').Value = '''0.0834'
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('D38').Value = '''3.16'
$ws.Range('E38').Value = '  -4.43%  '
$ws.Range('D39').Value = '''1.99'
$ws.Range('E39').Value = '  -7.08%  '
$ws.Range('D40').Value = '''17.78'
$ws.Range('E40').Value = '  -5.58%  '
$ws.Range('E41').Value = '  -1.76%  '
$ws.Range('D42').Value = '''125.32'
$ws.Range('E42').Value = '  -1.13%  '
$ws.Range('D43').Value = '''2.51'
$ws.Range('E43').Value = '  -4.80%  '
$ws.Range('E44').Value = '  -1.64%  '
$ws.Range('D45').Value = '''21.92'
$ws.Range('E45').Value = '  -7.80%  '
$ws.Range('D46').Value = '2.044.54'
$ws.Range('E46').Value = '  -1.86%  '
$ws.Range('D47').Value = '''2.33'
$ws.Range('E47').Value = '  +0.56%  '
$ws.Range('D48').Value = '''3.21'
$ws.Range('E48').Value = '  -5.55%  '
$ws.Range('D49').Value = '''5.73'
$ws.Range('E49').Value = '  +1.26%  '
$ws.Range('D50').Value = '''0.917'
$ws.Range('E50').Value = '  -4.83%  '
$ws.Range('D51').Value = '''8.93'
$ws.Range('E51').Value = '  -0.63%  '
